$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 172.33333
$ws.Range("I33").Value = 139.6
$ws.Range("J33").Value = 336
$ws.Range("K33").Value = 139.6
$ws.Range("L33").Value = 336
$ws.Range("M33").Value = 89.40000000000001
$ws.Range("N33").Value = -794
$ws.Range("H70").Value = 1750.2554
$ws.Range("I70").Value = 1665.5217
$ws.Range("K70").Value = 4996.5651
$ws.Range("M70").Value = -4726.5651
$ws.Range("H73").Value = 1750.2554
$ws.Range("I73").Value = 1665.5217
$ws.Range("K73").Value = 4996.5651
$ws.Range("M73").Value = -4060.5651
$ws.Range("H113").Value = 3234.923
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3783.7778
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3783.7778
$ws.Range("M113").Value = 1254
$ws.Range("N113").Value = -10291.7778
$ws.Range("H125").Value = 1249.75
$ws.Range("I125").Value = 499
$ws.Range("J125").Value = 1500
$ws.Range("K125").Value = 4491
$ws.Range("L125").Value = 13500
$ws.Range("M125").Value = -2031
$ws.Range("N125").Value = -18420
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25435.256
$ws.Range("I32").Value = 25435.256
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 25435.256
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -25148.256
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 9428.553
$ws.Range("I61").Value = 5620.136
$ws.Range("J61").Value = 14665.125
$ws.Range("K61").Value = 5620.136
$ws.Range("L61").Value = 14665.125
$ws.Range("M61").Value = -5408.136
$ws.Range("N61").Value = -15089.125
$ws.Range("H120").Value = 44766.332
$ws.Range("J120").Value = 44766.332
$ws.Range("L120").Value = 44766.332
$ws.Range("N120").Value = -54442.332
$ws.Range("H122").Value = 2963
$ws.Range("I122").Value = 2957.7144
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8873.143199999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6423.143199999999
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 11872.091
$ws.Range("I132").Value = 13731
$ws.Range("J132").Value = 3507
$ws.Range("K132").Value = 41193
$ws.Range("L132").Value = 10521
$ws.Range("M132").Value = -38663
$ws.Range("N132").Value = -15581
$ws.Range("H136").Value = 9428.553
$ws.Range("I136").Value = 5620.136
$ws.Range("J136").Value = 14665.125
$ws.Range("K136").Value = 16860.408
$ws.Range("L136").Value = 43995.375
$ws.Range("M136").Value = -14310.408
$ws.Range("N136").Value = -49095.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 59421.832
$ws.Range("I134").Value = 4984.875
$ws.Range("J134").Value = 102971.4
$ws.Range("K134").Value = 14954.625
$ws.Range("L134").Value = 308914.2
$ws.Range("M134").Value = -12419.625
$ws.Range("N134").Value = -313984.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3503143.8
$ws.Range("I58").Value = 5052920.5
$ws.Range("K58").Value = 5052920.5
$ws.Range("M58").Value = -5052717.5
$ws.Range("H68").Value = 40196.668
$ws.Range("J68").Value = 40196.668
$ws.Range("L68").Value = 40196.668
$ws.Range("N68").Value = -41694.668
$ws.Range("H71").Value = 40196.668
$ws.Range("J71").Value = 40196.668
$ws.Range("L71").Value = 120590.004
$ws.Range("N71").Value = -128078.004
$ws.Range("H99").Value = 3520
$ws.Range("I99").Value = 3150
$ws.Range("J99").Value = 3766.6667
$ws.Range("K99").Value = 3150
$ws.Range("L99").Value = 3766.6667
$ws.Range("M99").Value = -1652
$ws.Range("N99").Value = -6762.6667
$ws.Range("H126").Value = 3520
$ws.Range("I126").Value = 3150
$ws.Range("J126").Value = 3766.6667
$ws.Range("K126").Value = 9450
$ws.Range("L126").Value = 11300.0001
$ws.Range("M126").Value = -6980
$ws.Range("N126").Value = -16240.0001
$ws.Range("H136").Value = 3503143.8
$ws.Range("I136").Value = 5052920.5
$ws.Range("K136").Value = 15158761.5
$ws.Range("M136").Value = -15156211.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1050.45
$ws.Range("I131").Value = 1143.2222
$ws.Range("K131").Value = 3429.6666
$ws.Range("M131").Value = 1610.3334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 2500
$ws.Range("J92").Value = 2500
$ws.Range("L92").Value = 2500
$ws.Range("N92").Value = -6244
$ws.Range("H122").Value = 22999.6
$ws.Range("I122").Value = 27249.5
$ws.Range("K122").Value = 81748.5
$ws.Range("M122").Value = -79298.5
$ws.Range("H126").Value = 3257.1428
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 4200
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 12600
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -17540
$ws.Range("H132").Value = 22374.092
$ws.Range("I132").Value = 16762.625
$ws.Range("J132").Value = 37338
$ws.Range("K132").Value = 50287.875
$ws.Range("L132").Value = 112014
$ws.Range("M132").Value = -47757.875
$ws.Range("N132").Value = -117074
$ws.Range("H139").Value = 43999.8
$ws.Range("J139").Value = 43999.8
$ws.Range("L139").Value = 43999.8
$ws.Range("N139").Value = -54279.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 686
$ws.Range("I22").Value = 604
$ws.Range("J22").Value = 822.6667
$ws.Range("K22").Value = 604
$ws.Range("L22").Value = 822.6667
$ws.Range("M22").Value = -309
$ws.Range("N22").Value = -1412.6667
$ws.Range("H27").Value = 686
$ws.Range("I27").Value = 604
$ws.Range("J27").Value = 822.6667
$ws.Range("K27").Value = 604
$ws.Range("L27").Value = 822.6667
$ws.Range("M27").Value = -497
$ws.Range("N27").Value = -1036.6667
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2864
$ws.Range("H132").Value = 5929.1
$ws.Range("I132").Value = 6252
$ws.Range("J132").Value = 5713.8335
$ws.Range("K132").Value = 18756
$ws.Range("L132").Value = 17141.5005
$ws.Range("M132").Value = -16226
$ws.Range("N132").Value = -22201.5005
$ws.Range("H136").Value = 5851.0938
$ws.Range("I136").Value = 3310.3845
$ws.Range("J136").Value = 7589.4736
$ws.Range("K136").Value = 9931.1535
$ws.Range("L136").Value = 22768.4208
$ws.Range("M136").Value = -7381.1535
$ws.Range("N136").Value = -27868.4208
$ws.Range("H140").Value = 72103.75
$ws.Range("J140").Value = 72103.75
$ws.Range("L140").Value = 72103.75
$ws.Range("N140").Value = -82463.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3801.1
$ws.Range("J62").Value = 3502.25
$ws.Range("L62").Value = 3502.25
$ws.Range("N62").Value = -4750.25
$ws.Range("H65").Value = 3801.1
$ws.Range("J65").Value = 3502.25
$ws.Range("L65").Value = 17511.25
$ws.Range("N65").Value = -23751.25
$ws.Range("H124").Value = 21371.6
$ws.Range("J124").Value = 21371.6
$ws.Range("L124").Value = 21371.6
$ws.Range("N124").Value = -31191.6
$ws.Range("H126").Value = 1655.9584
$ws.Range("I126").Value = 1822
$ws.Range("J126").Value = 1423.5
$ws.Range("K126").Value = 5466
$ws.Range("L126").Value = 4270.5
$ws.Range("M126").Value = -2996
$ws.Range("N126").Value = -9210.5
$ws.Range("H132").Value = 7834.5
$ws.Range("I132").Value = 9001.333000000001
$ws.Range("K132").Value = 27003.999
$ws.Range("M132").Value = -24473.999
$ws.Range("H136").Value = 5437.9346
$ws.Range("I136").Value = 2373.2068
$ws.Range("K136").Value = 7119.6204
$ws.Range("M136").Value = -4569.6204
